$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "US English"
$ws.Range("C17").Value = 24.5
$ws.Range("D17").Value = 18.9
$ws.Range("F17").Value = 6.6
$ws.Range("G17").Value = 6.5
$ws.Range("I17").Value = 9.6
$ws.Range("J17").Value = 90.4

$ws.Range("A18").Value = "UK English"
$ws.Range("C18").Value = 26.9
$ws.Range("D18").Value = 11.9
$ws.Range("F18").Value = 8.6
$ws.Range("G18").Value = 4.2
$ws.Range("I18").Value = 12.5
$ws.Range("J18").Value = 87.5

$styleRange = $ws.Range("A17:D18")
$styleRange.Font.Name = ".AppleSystemUIFont"
$styleRange.Font.Size = 12
$styleRange.Font.Color = 4539717

$ws.Range("B17:B18").Clear()

$ws.Range("G22").Select()
